$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 108.588234
$ws.Range("J33").Value = 166.4
$ws.Range("L33").Value = 166.4
$ws.Range("N33").Value = -624.4
$ws.Range("H43").Value = 2038
$ws.Range("I43").Value = 1319.5
$ws.Range("J43").Value = 2397.25
$ws.Range("K43").Value = 1319.5
$ws.Range("L43").Value = 2397.25
$ws.Range("M43").Value = -1250.5
$ws.Range("N43").Value = -2535.25
$ws.Range("H80").Value = 659.1818
$ws.Range("I80").Value = 723.8
$ws.Range("K80").Value = 2171.4
$ws.Range("M80").Value = -1173.4
$ws.Range("H83").Value = 659.1818
$ws.Range("I83").Value = 723.8
$ws.Range("K83").Value = 6514.2
$ws.Range("M83").Value = -1522.2
$ws.Range("H95").Value = 38332.332
$ws.Range("J95").Value = 38332.332
$ws.Range("L95").Value = 38332.332
$ws.Range("N95").Value = -43824.332
$ws.Range("H97").Value = 1054.5
$ws.Range("J97").Value = 1054.5
$ws.Range("L97").Value = 3163.5
$ws.Range("N97").Value = -4155.5
$ws.Range("H112").Value = 3632.5
$ws.Range("I112").Value = 5000
$ws.Range("J112").Value = 3359
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 10077
$ws.Range("M112").Value = -13892
$ws.Range("N112").Value = -12293
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H138").Value = 2287.1538
$ws.Range("J138").Value = 15000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3032874
$ws.Range("I32").Value = 586.2963
$ws.Range("J32").Value = 16678169
$ws.Range("K32").Value = 586.2963
$ws.Range("L32").Value = 16678169
$ws.Range("M32").Value = -299.2963
$ws.Range("N32").Value = -16678743
$ws.Range("H44").Value = 12783.375
$ws.Range("J44").Value = 12783.375
$ws.Range("L44").Value = 12783.375
$ws.Range("N44").Value = -13759.375
$ws.Range("H55").Value = 99999
$ws.Range("J55").Value = 99999
$ws.Range("L55").Value = 99999
$ws.Range("N55").Value = -100629
$ws.Range("H102").Value = 6253066.5
$ws.Range("I102").Value = 8929381
$ws.Range("K102").Value = 8929381
$ws.Range("M102").Value = -8927759
$ws.Range("H110").Value = 90912000
$ws.Range("J110").Value = 2623.625
$ws.Range("L110").Value = 2623.625
$ws.Range("N110").Value = -6713.625
$ws.Range("H122").Value = 2044.9166
$ws.Range("I122").Value = 1388.6666
$ws.Range("K122").Value = 4165.9998
$ws.Range("M122").Value = -1715.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 726.6667
$ws.Range("I20").Value = 590
$ws.Range("K20").Value = 590
$ws.Range("M20").Value = -343
$ws.Range("H25").Value = 5983.3335
$ws.Range("I25").Value = 2300
$ws.Range("K25").Value = 2300
$ws.Range("M25").Value = -2065
$ws.Range("H54").Value = 13500
$ws.Range("I54").Value = 13500
$ws.Range("K54").Value = 13500
$ws.Range("M54").Value = -13016
$ws.Range("H86").Value = 2734.65
$ws.Range("I86").Value = 1384.5
$ws.Range("K86").Value = 1384.5
$ws.Range("M86").Value = -261.5
$ws.Range("H89").Value = 2734.65
$ws.Range("I89").Value = 1384.5
$ws.Range("K89").Value = 6922.5
$ws.Range("M89").Value = -1306.5
$ws.Range("H94").Value = 700
$ws.Range("I94").Value = 700
$ws.Range("K94").Value = 700
$ws.Range("M94").Value = -249
$ws.Range("H99").Value = 100001016
$ws.Range("I99").Value = 125001010
$ws.Range("J99").Value = 1082.5
$ws.Range("K99").Value = 125001010
$ws.Range("L99").Value = 1082.5
$ws.Range("M99").Value = -124999512
$ws.Range("N99").Value = -4078.5
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988
$ws.Range("H107").Value = 26321106
$ws.Range("I107").Value = 62502124
$ws.Range("K107").Value = 62502124
$ws.Range("M107").Value = -62500204
$ws.Range("H134").Value = 4309.6523
$ws.Range("I134").Value = 1379.1111
$ws.Range("J134").Value = 14859.6
$ws.Range("K134").Value = 4137.3333
$ws.Range("L134").Value = 44578.8
$ws.Range("M134").Value = -1602.3333
$ws.Range("N134").Value = -49648.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3246.5
$ws.Range("I31").Value = 1893
$ws.Range("K31").Value = 1893
$ws.Range("M31").Value = -1598
$ws.Range("H34").Value = 3246.5
$ws.Range("I34").Value = 1893
$ws.Range("K34").Value = 1893
$ws.Range("M34").Value = -1691
$ws.Range("H58").Value = 2682
$ws.Range("I58").Value = 1579.6364
$ws.Range("K58").Value = 1579.6364
$ws.Range("M58").Value = -1376.6364
$ws.Range("H122").Value = 637
$ws.Range("I122").Value = 512
$ws.Range("J122").Value = 699.5
$ws.Range("K122").Value = 1536
$ws.Range("L122").Value = 2098.5
$ws.Range("M122").Value = 914
$ws.Range("N122").Value = -6998.5
$ws.Range("H123").Value = 44892
$ws.Range("J123").Value = 44892
$ws.Range("L123").Value = 44892
$ws.Range("N123").Value = -54692
$ws.Range("H136").Value = 2682
$ws.Range("I136").Value = 1579.6364
$ws.Range("K136").Value = 4738.9092
$ws.Range("M136").Value = -2188.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2337.5
$ws.Range("I57").Value = 1675
$ws.Range("K57").Value = 5025
$ws.Range("M57").Value = -4466
$ws.Range("H104").Value = 9049.9
$ws.Range("J104").Value = 9944.333000000001
$ws.Range("L104").Value = 29832.999
$ws.Range("N104").Value = -35074.999
$ws.Range("H131").Value = 3177.8
$ws.Range("I131").Value = 2000
$ws.Range("K131").Value = 6000
$ws.Range("M131").Value = -960
$ws.Range("H136").Value = 10333
$ws.Range("I136").Value = 8999.5
$ws.Range("K136").Value = 26998.5
$ws.Range("M136").Value = -21898.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6250
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 7500
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 7500
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -9496
$ws.Range("H83").Value = 6250
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 7500
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 37500
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -47484
$ws.Range("H102").Value = 2131.65
$ws.Range("I102").Value = 1980.6842
$ws.Range("K102").Value = 1980.6842
$ws.Range("M102").Value = -358.6841999999999
$ws.Range("H122").Value = 1936.6364
$ws.Range("I122").Value = 1197.1428
$ws.Range("J122").Value = 3230.75
$ws.Range("K122").Value = 3591.4284
$ws.Range("L122").Value = 9692.25
$ws.Range("M122").Value = -1141.4284
$ws.Range("N122").Value = -14592.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2310
$ws.Range("I22").Value = 1775
$ws.Range("K22").Value = 1775
$ws.Range("M22").Value = -1480
$ws.Range("H27").Value = 2310
$ws.Range("I27").Value = 1775
$ws.Range("K27").Value = 1775
$ws.Range("M27").Value = -1668
$ws.Range("I61").Value = 125001230
$ws.Range("K61").Value = 125001230
$ws.Range("M61").Value = -125001028
$ws.Range("H93").Value = 3948.6667
$ws.Range("I93").Value = 3421.5
$ws.Range("K93").Value = 3421.5
$ws.Range("M93").Value = -2173.5
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180
$ws.Range("I113").Value = 125001230
$ws.Range("K113").Value = 125001230
$ws.Range("M113").Value = -124999060
$ws.Range("H132").Value = 2403.45
$ws.Range("I132").Value = 1269.9375
$ws.Range("K132").Value = 3809.8125
$ws.Range("M132").Value = -1279.8125
